$d = $word.ActiveDocument

$d.Content.Find.Execute("176÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "694÷8=", 2) | Out-Null
$d.Content.Find.Execute("538÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "269÷7=", 2) | Out-Null
$d.Content.Find.Execute("566÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "537÷6=", 2) | Out-Null
$d.Content.Find.Execute("429÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "930÷8=", 2) | Out-Null
$d.Content.Find.Execute("221÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "900÷5=", 2) | Out-Null
$d.Content.Find.Execute("182÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "872÷4=", 2) | Out-Null
$d.Content.Find.Execute("797÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "672÷2=", 2) | Out-Null
$d.Content.Find.Execute("945÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "453÷3=", 2) | Out-Null
$d.Content.Find.Execute("589÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "629÷6=", 2) | Out-Null
$d.Content.Find.Execute("564÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "959÷6=", 2) | Out-Null
$d.Content.Find.Execute("818÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "640÷5=", 2) | Out-Null
$d.Content.Find.Execute("993÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "445÷4=", 2) | Out-Null
$d.Content.Find.Execute("577÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "879÷4=", 2) | Out-Null
$d.Content.Find.Execute("138÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "215÷6=", 2) | Out-Null
$d.Content.Find.Execute("736÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "503÷3=", 2) | Out-Null
$d.Content.Find.Execute("741÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "482÷8=", 2) | Out-Null
$d.Content.Find.Execute("879÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "167÷2=", 2) | Out-Null
$d.Content.Find.Execute("596÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "440÷9=", 2) | Out-Null
$d.Content.Find.Execute("318÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "191÷8=", 2) | Out-Null
$d.Content.Find.Execute("698÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "521÷6=", 2) | Out-Null
$d.Content.Find.Execute("439÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "208÷7=", 2) | Out-Null
$d.Content.Find.Execute("524÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "270÷4=", 2) | Out-Null
$d.Content.Find.Execute("311÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "109÷8=", 2) | Out-Null
$d.Content.Find.Execute("888÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "292÷7=", 2) | Out-Null
$d.Content.Find.Execute("631÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "578÷8=", 2) | Out-Null
